# Updates the "cryptos" price/volume snapshot on Sheet1 (rows 2-51).
# For D-column cells whose new text looks like a plain number (e.g. "241.12"),
# the cell is briefly switched to Text format before the assignment so Excel
# stores it as text (matching the source data's inlineStr string cells)
# instead of silently converting it to a numeric value; the style is then
# reset to "Normal" so no stray number-format style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '96.654.88'
$ws.Range("E2").Value = '  +1.85%  '

$ws.Range("D3").Value = '3.571.73'
$ws.Range("E3").Value = '  -0.98%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '653.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.54%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.65'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +13.97%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.412'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.10%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.07'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.53%  '

$ws.Range("B10").Value = 'USDC'
$ws.Range("C10").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.00%  '

$ws.Range("D11").Value = '3.570.73'
$ws.Range("E11").Value = '  -0.98%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.55'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.81%  '

$ws.Range("E13").Value = '  +1.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.17%  '

$ws.Range("D15").Value = '4.235.17'
$ws.Range("E15").Value = '  -1.37%  '

$ws.Range("D16").Value = '96.457.69'
$ws.Range("E16").Value = '  +1.58%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000261'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.29%  '

$ws.Range("D18").Value = '3.543.76'
$ws.Range("E18").Value = '  -1.70%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.532'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +11.51%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '508.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.74%  '

$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.40'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.35%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000200'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '96.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.61%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.32%  '

$ws.Range("D29").Value = '3.763.28'
$ws.Range("E29").Value = '  -0.88%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.156'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +13.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.02'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.46'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.07%  '

$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.184'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.77%  '

$ws.Range("E35").Value = '  +0.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.44'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '627.08'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.82%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.568'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.64'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.56%  '

$ws.Range("E41").Value = '  +0.06%  '

$ws.Range("E42").Value = '  +1.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.906'
$ws.Range("D43").Style = "Normal"

$ws.Range("E44").Value = '  +6.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0427'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.30'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.54'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.06'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.94%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.52'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.30'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.63%  '
